# Auto-generated script to apply price/profit recalculation updates
# to the Ultima_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 633.1667
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 679.8
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 679.8
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -905.8
$ws.Range("H113").Value = 3332.5
$ws.Range("I113").Value = 2665.3125
$ws.Range("J113").Value = 4400
$ws.Range("K113").Value = 2665.3125
$ws.Range("L113").Value = 4400
$ws.Range("M113").Value = 588.6875
$ws.Range("N113").Value = -10908
$ws.Range("H130").Value = 70780
$ws.Range("J130").Value = 70780
$ws.Range("L130").Value = 70780
$ws.Range("N130").Value = -80820
$ws.Range("H133").Value = 34825.715
$ws.Range("J133").Value = 34825.715
$ws.Range("L133").Value = 34825.715
$ws.Range("N133").Value = -44945.715
$ws.Range("H134").Value = 53000
$ws.Range("J134").Value = 53000
$ws.Range("L134").Value = 53000
$ws.Range("N134").Value = -63140
$ws.Range("H136").Value = 64678.57
$ws.Range("J136").Value = 64678.57
$ws.Range("L136").Value = 64678.57
$ws.Range("N136").Value = -74878.57000000001
$ws.Range("H137").Value = 8482
$ws.Range("I137").Value = 1471.7142
$ws.Range("K137").Value = 4415.142599999999
$ws.Range("M137").Value = -1865.142599999999
$ws.Range("H139").Value = 45179.855
$ws.Range("J139").Value = 45179.855
$ws.Range("L139").Value = 45179.855
$ws.Range("N139").Value = -55459.855
$ws.Range("H140").Value = 50043.75
$ws.Range("I140").Value = 12000
$ws.Range("J140").Value = 55478.57
$ws.Range("K140").Value = 12000
$ws.Range("L140").Value = 55478.57
$ws.Range("M140").Value = -6820
$ws.Range("N140").Value = -65838.57000000001
$ws.Range("H141").Value = 1229.0385
$ws.Range("I141").Value = 1019.7826
$ws.Range("J141").Value = 2833.3333
$ws.Range("K141").Value = 3059.3478
$ws.Range("L141").Value = 8499.999899999999
$ws.Range("M141").Value = 2120.6522
$ws.Range("N141").Value = -18859.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1294.6364
$ws.Range("I74").Value = 1346.421
$ws.Range("J74").Value = 966.6667
$ws.Range("K74").Value = 1346.421
$ws.Range("L74").Value = 966.6667
$ws.Range("M74").Value = -472.421
$ws.Range("N74").Value = -2714.6667
$ws.Range("H77").Value = 1294.6364
$ws.Range("I77").Value = 1346.421
$ws.Range("J77").Value = 966.6667
$ws.Range("K77").Value = 6732.105
$ws.Range("L77").Value = 4833.3335
$ws.Range("M77").Value = -2364.105
$ws.Range("N77").Value = -13569.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 989.75
$ws.Range("I94").Value = 729.6667
$ws.Range("K94").Value = 729.6667
$ws.Range("M94").Value = -278.6667
$ws.Range("H105").Value = 3905.1025
$ws.Range("I105").Value = 2141.111
$ws.Range("J105").Value = 4434.3
$ws.Range("K105").Value = 2141.111
$ws.Range("L105").Value = 4434.3
$ws.Range("M105").Value = -394.1109999999999
$ws.Range("N105").Value = -7928.3
$ws.Range("H134").Value = 3583.3447
$ws.Range("I134").Value = 2246.4707
$ws.Range("J134").Value = 5477.25
$ws.Range("K134").Value = 6739.4121
$ws.Range("L134").Value = 16431.75
$ws.Range("M134").Value = -4204.4121
$ws.Range("N134").Value = -21501.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3760.3704
$ws.Range("I31").Value = 2271.4285
$ws.Range("K31").Value = 2271.4285
$ws.Range("M31").Value = -1976.4285
$ws.Range("H34").Value = 3760.3704
$ws.Range("I34").Value = 2271.4285
$ws.Range("K34").Value = 2271.4285
$ws.Range("M34").Value = -2069.4285
$ws.Range("H41").Value = 22666.666
$ws.Range("H50").Value = 9389.200000000001
$ws.Range("J50").Value = 9389.200000000001
$ws.Range("L50").Value = 9389.200000000001
$ws.Range("N50").Value = -10639.2
$ws.Range("H51").Value = 9295.4
$ws.Range("J51").Value = 9295.4
$ws.Range("L51").Value = 9295.4
$ws.Range("N51").Value = -10767.4
$ws.Range("H59").Value = 15651
$ws.Range("J59").Value = 15651
$ws.Range("L59").Value = 15651
$ws.Range("N59").Value = -17941
$ws.Range("H60").Value = 8040.4287
$ws.Range("J60").Value = 8256.6
$ws.Range("L60").Value = 8256.6
$ws.Range("N60").Value = -9278.6
$ws.Range("H61").Value = 9295.4
$ws.Range("J61").Value = 9295.4
$ws.Range("L61").Value = 9295.4
$ws.Range("N61").Value = -9991.4
$ws.Range("H68").Value = 17249.666
$ws.Range("J68").Value = 17249.666
$ws.Range("L68").Value = 17249.666
$ws.Range("N68").Value = -18747.666
$ws.Range("H71").Value = 17249.666
$ws.Range("J71").Value = 17249.666
$ws.Range("L71").Value = 51748.99800000001
$ws.Range("N71").Value = -59236.99800000001
$ws.Range("H74").Value = 13702.25
$ws.Range("J74").Value = 13702.25
$ws.Range("L74").Value = 13702.25
$ws.Range("N74").Value = -15450.25
$ws.Range("H77").Value = 13702.25
$ws.Range("J77").Value = 13702.25
$ws.Range("L77").Value = 41106.75
$ws.Range("N77").Value = -49842.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 793.05554
$ws.Range("I5").Value = 317.5
$ws.Range("J5").Value = 1387.5
$ws.Range("K5").Value = 952.5
$ws.Range("L5").Value = 4162.5
$ws.Range("M5").Value = -840.5
$ws.Range("N5").Value = -4386.5
$ws.Range("H131").Value = 1775.2253
$ws.Range("I131").Value = 2711.9473
$ws.Range("J131").Value = 1432.9615
$ws.Range("K131").Value = 8135.841899999999
$ws.Range("L131").Value = 4298.8845
$ws.Range("M131").Value = -3095.841899999999
$ws.Range("N131").Value = -14378.8845
$ws.Range("H135").Value = 793.05554
$ws.Range("I135").Value = 317.5
$ws.Range("J135").Value = 1387.5
$ws.Range("K135").Value = 2857.5
$ws.Range("L135").Value = 12487.5
$ws.Range("M135").Value = -322.5
$ws.Range("N135").Value = -17557.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10907.821
$ws.Range("I70").Value = 12815.5
$ws.Range("J70").Value = 3913
$ws.Range("K70").Value = 12815.5
$ws.Range("L70").Value = 3913
$ws.Range("M70").Value = -12545.5
$ws.Range("N70").Value = -4453
$ws.Range("H73").Value = 10907.821
$ws.Range("I73").Value = 12815.5
$ws.Range("J73").Value = 3913
$ws.Range("K73").Value = 12815.5
$ws.Range("L73").Value = 3913
$ws.Range("M73").Value = -11879.5
$ws.Range("N73").Value = -5785
$ws.Range("H97").Value = 1980.2727
$ws.Range("I97").Value = 1433.3334
$ws.Range("J97").Value = 2185.375
$ws.Range("K97").Value = 1433.3334
$ws.Range("L97").Value = 2185.375
$ws.Range("M97").Value = -937.3334
$ws.Range("N97").Value = -3177.375
$ws.Range("H113").Value = 501405.5
$ws.Range("I113").Value = 1000011
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1000011
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -997841
$ws.Range("N113").Value = -7140
$ws.Range("H122").Value = 2883.7058
$ws.Range("I122").Value = 2486.0908
$ws.Range("J122").Value = 3612.6667
$ws.Range("K122").Value = 7458.2724
$ws.Range("L122").Value = 10838.0001
$ws.Range("M122").Value = -5008.2724
$ws.Range("N122").Value = -15738.0001
$ws.Range("H132").Value = 5869.483
$ws.Range("I132").Value = 6284.1665
$ws.Range("J132").Value = 3879
$ws.Range("K132").Value = 18852.4995
$ws.Range("L132").Value = 11637
$ws.Range("M132").Value = -16322.4995
$ws.Range("N132").Value = -16697

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2790.1538
$ws.Range("J40").Value = 1542.6786
$ws.Range("L40").Value = 1542.6786
$ws.Range("N40").Value = -1814.6786
$ws.Range("H46").Value = 1143.4615
$ws.Range("I46").Value = 1098
$ws.Range("J46").Value = 1182.4286
$ws.Range("K46").Value = 1098
$ws.Range("L46").Value = 1182.4286
$ws.Range("M46").Value = -910
$ws.Range("N46").Value = -1558.4286
$ws.Range("H132").Value = 8338798
$ws.Range("I132").Value = 3762.5652
$ws.Range("J132").Value = 35725344
$ws.Range("K132").Value = 11287.6956
$ws.Range("L132").Value = 107176032
$ws.Range("M132").Value = -8757.695599999999
$ws.Range("N132").Value = -107181092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 544.7778
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 420.6
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 1261.8
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -5601.8
